$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# A new source file (9ad15719-6ce1-48cc-8569-036f14eacdc0.md) was handed off,
# in addition to the previously-tracked a6289566-8408-4fec-9a8c-04180367b81e.md.
# A row for it is inserted above the existing a6289566 row on every sheet,
# pushing the existing rows down by one.
# ---------------------------------------------------------------------------

# =================================== Overview ===================================
$wsOv = $wb.Worksheets.Item("Overview")

# Drop the existing hyperlinks before shifting rows around - this engine does not
# auto-shift hyperlink anchors when rows are inserted, so they are re-created
# from scratch afterwards, in final top-to-bottom order.
$wsOv.Range("A2").Hyperlinks.Delete()
$wsOv.Range("A3").Hyperlinks.Delete()

$wsOv.Rows.Item(2).Insert()

$wsOv.Range("A2").Value = "9ad15719-6ce1-48cc-8569-036f14eacdc0.md"
$wsOv.Range("B2").Value = "Ready for handoff"
$wsOv.Range("C2").Value = "Ready for handoff"

$wsOv.Hyperlinks.Add($wsOv.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/29f8223291f41775c48fc5f81140839eeff1258d/e2e/9ad15719-6ce1-48cc-8569-036f14eacdc0.md", "", "", "9ad15719-6ce1-48cc-8569-036f14eacdc0.md") | Out-Null
$wsOv.Hyperlinks.Add($wsOv.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/29f8223291f41775c48fc5f81140839eeff1258d/e2e/a6289566-8408-4fec-9a8c-04180367b81e.md", "", "", "a6289566-8408-4fec-9a8c-04180367b81e.md") | Out-Null
$wsOv.Hyperlinks.Add($wsOv.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/29f8223291f41775c48fc5f81140839eeff1258d/.localization-config", "", "", ".localization-config") | Out-Null

# =================================== zh-cn =======================================
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("C2").Hyperlinks.Delete()
$wsZh.Range("A3").Hyperlinks.Delete()

$wsZh.Rows.Item(2).Insert()

$wsZh.Range("A2").Value = "9ad15719-6ce1-48cc-8569-036f14eacdc0.md"
$wsZh.Range("B2").Value = "Ready for handoff"
$wsZh.Range("C2").Value = "9ad15719-6ce1-48cc-8569-036f14eacdc0.79be59d77f61d35d36ff06509bb7aaec296dacce.zh-cn.xlf"
$wsZh.Range("D2").Value = "2016-03-09 20:33:08"
$wsZh.Range("G2").Value = "0001-01-01 00:00:00"
$wsZh.Range("H2").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/29f8223291f41775c48fc5f81140839eeff1258d/e2e/9ad15719-6ce1-48cc-8569-036f14eacdc0.md", "", "", "9ad15719-6ce1-48cc-8569-036f14eacdc0.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/79be59d77f61d35d36ff06509bb7aaec296dacce/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9ad15719-6ce1-48cc-8569-036f14eacdc0.79be59d77f61d35d36ff06509bb7aaec296dacce.zh-cn.xlf", "", "", "9ad15719-6ce1-48cc-8569-036f14eacdc0.79be59d77f61d35d36ff06509bb7aaec296dacce.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/29f8223291f41775c48fc5f81140839eeff1258d/e2e/a6289566-8408-4fec-9a8c-04180367b81e.md", "", "", "a6289566-8408-4fec-9a8c-04180367b81e.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/16eb241e7716f10c06e8c8193ec79918d1a355bd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a6289566-8408-4fec-9a8c-04180367b81e.12bacf9ab37516007f665f582b00427400306d74.zh-cn.xlf", "", "", "a6289566-8408-4fec-9a8c-04180367b81e.12bacf9ab37516007f665f582b00427400306d74.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/29f8223291f41775c48fc5f81140839eeff1258d/.localization-config", "", "", ".localization-config") | Out-Null

# =================================== de-de =======================================
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("C2").Hyperlinks.Delete()
$wsDe.Range("A3").Hyperlinks.Delete()

$wsDe.Rows.Item(2).Insert()

$wsDe.Range("A2").Value = "9ad15719-6ce1-48cc-8569-036f14eacdc0.md"
$wsDe.Range("B2").Value = "Ready for handoff"
$wsDe.Range("C2").Value = "9ad15719-6ce1-48cc-8569-036f14eacdc0.79be59d77f61d35d36ff06509bb7aaec296dacce.de-de.xlf"
$wsDe.Range("D2").Value = "2016-03-09 20:33:13"
$wsDe.Range("G2").Value = "0001-01-01 00:00:00"
$wsDe.Range("H2").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/29f8223291f41775c48fc5f81140839eeff1258d/e2e/9ad15719-6ce1-48cc-8569-036f14eacdc0.md", "", "", "9ad15719-6ce1-48cc-8569-036f14eacdc0.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/79be59d77f61d35d36ff06509bb7aaec296dacce/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9ad15719-6ce1-48cc-8569-036f14eacdc0.79be59d77f61d35d36ff06509bb7aaec296dacce.de-de.xlf", "", "", "9ad15719-6ce1-48cc-8569-036f14eacdc0.79be59d77f61d35d36ff06509bb7aaec296dacce.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/29f8223291f41775c48fc5f81140839eeff1258d/e2e/a6289566-8408-4fec-9a8c-04180367b81e.md", "", "", "a6289566-8408-4fec-9a8c-04180367b81e.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5eabcd2777cdee4bb005f0839057034f59542b7e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a6289566-8408-4fec-9a8c-04180367b81e.12bacf9ab37516007f665f582b00427400306d74.de-de.xlf", "", "", "a6289566-8408-4fec-9a8c-04180367b81e.12bacf9ab37516007f665f582b00427400306d74.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/29f8223291f41775c48fc5f81140839eeff1258d/.localization-config", "", "", ".localization-config") | Out-Null
